$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $value) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "66.233.09"
Set-TextValue "E2" "  +0.73%  "
Set-TextValue "D3" "3.322.72"
Set-TextValue "E3" "  +0.87%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.07%  "
Set-TextValue "D5" "564.41"
Set-TextValue "D6" "186.24"
Set-TextValue "E6" "  +1.62%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "D8" "3.317.74"
Set-TextValue "E8" "  +0.87%  "
Set-TextValue "E9" "  -1.94%  "
Set-TextValue "E10" "  -4.80%  "
Set-TextValue "D11" "0.575"
Set-TextValue "E11" "  -1.62%  "
Set-TextValue "D12" "46.16"
Set-TextValue "E12" "  -2.46%  "
Set-TextValue "D13" "0.0000266"
Set-TextValue "E13" "  -0.07%  "
Set-TextValue "D14" "3.852.83"
Set-TextValue "E14" "  +0.88%  "
Set-TextValue "D15" "8.46"
Set-TextValue "E15" "  -2.06%  "
Set-TextValue "D16" "595.30"
Set-TextValue "E16" "  -7.54%  "
Set-TextValue "D17" "66.228.33"
Set-TextValue "E17" "  +0.77%  "
Set-TextValue "E18" "  +0.38%  "
Set-TextValue "D19" "3.320.27"
Set-TextValue "E19" "  +0.74%  "
Set-TextValue "D20" "17.73"
Set-TextValue "E20" "  -2.14%  "
Set-TextValue "D21" "10.94"
Set-TextValue "E21" "  -3.97%  "
Set-TextValue "D22" "0.898"
Set-TextValue "E22" "  -0.79%  "
Set-TextValue "D23" "17.96"
Set-TextValue "E23" "  -1.34%  "
Set-TextValue "D24" "4.98"
Set-TextValue "E24" "  +1.57%  "
Set-TextValue "D25" "98.69"
Set-TextValue "E25" "  -8.22%  "
Set-TextValue "D26" "3.99"
Set-TextValue "E26" "  +0.25%  "
Set-TextValue "D27" "2.71"
Set-TextValue "E27" "  +1.07%  "
Set-TextValue "D28" "9.41"
Set-TextValue "E28" "  -1.55%  "
Set-TextValue "D29" "8.46"
Set-TextValue "E29" "  -2.43%  "
Set-TextValue "D30" "30.71"
Set-TextValue "E30" "  +1.56%  "
Set-TextValue "E31" "  +5.54%  "
Set-TextValue "B32" "Bittensor"
Set-TextValue "C32" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D32" "563.15"
Set-TextValue "E32" "  +8.02%  "
Set-TextValue "B33" "dogwifhat"
Set-TextValue "C33" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D33" "3.67"
Set-TextValue "E33" "  -6.12%  "
Set-TextValue "D34" "10.86"
Set-TextValue "E34" "  -1.91%  "
Set-TextValue "D35" "3.809.40"
Set-TextValue "E35" "  +0.25%  "
Set-TextValue "E36" "  -1.28%  "
Set-TextValue "E37" "  +0.05%  "
Set-TextValue "D38" "55.87"
Set-TextValue "E38" "  -2.78%  "
Set-TextValue "D39" "33.30"
Set-TextValue "E39" "  +1.08%  "
Set-TextValue "E40" "  -2.14%  "
Set-TextValue "D41" "0.0₃0689"
Set-TextValue "E41" "  -6.10%  "
Set-TextValue "E42" "  -6.42%  "
Set-TextValue "D43" "3.39"
Set-TextValue "E43" "  +4.08%  "
Set-TextValue "D44" "2.60"
Set-TextValue "E44" "  -4.36%  "
Set-TextValue "D45" "0.333"
Set-TextValue "E45" "  -1.41%  "
Set-TextValue "E46" "  -0.82%  "
Set-TextValue "D47" "3.03"
Set-TextValue "E47" "  -9.62%  "
Set-TextValue "E48" "  -2.38%  "
Set-TextValue "E49" "  +0.09%  "
Set-TextValue "D50" "2.53"
Set-TextValue "E50" "  -2.90%  "
Set-TextValue "D51" "129.26"
Set-TextValue "E51" "  +6.11%  "
